# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
# Update the account-statement (Estado de Cuenta) worksheet:
#   - correct the mora period for the existing worker (2508 -> 2509)
#   - add a second worker row (ERNESTO JIMENEZ FRANCO) for the same period
#   - refresh the summary total (Valor Mora) and worker count accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row right below the existing worker row (row 16). This pushes
# the footer rows (and their merged ranges) down by one, exactly like a user
# inserting a row in the middle of the table.
$ws.Rows("17:17").Insert()

# Copy the formatting (borders, fills, fonts, number formats) from the
# existing worker row onto the newly inserted row so it matches the table.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fix the mora period on the existing worker row (2508 -> 2509).
$ws.Range("E16").Value = "2509"

# Populate the new worker row with the second debtor's data.
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73141106"
$ws.Range("D17").Value = "ERNESTO JIMENEZ FRANCO"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Refresh the header summary: total mora value and worker count now reflect
# both workers (56940 + 56940 = 113880; 1 -> 2 workers).
$ws.Range("E11").Value = 113880
$ws.Range("C13").Value = 2
